$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "Objetivos:" row (row 10) B/C text
$objetivosTxt = "Proporcionar conhecimentos e técnicas necessárias para o entendimento do funcionamento de softwares de apoio a execução de desenho técnico"
$ws.Cells.Item(10, 2).Value = $objetivosTxt
$ws.Cells.Item(10, 3).Value = $objetivosTxt

# 2) Insert a new row before row 13 ("Docentes responsáveis:" data row),
#    shifting the old rows 13-21 down to 14-22.
$ws.Rows.Item(13).Insert()

# Copy the formatting from row 2 (B2:C2, which carries styles 2/3 with no
# value in column A) onto the newly inserted row 13 so the new row matches
# the target look (normal black B, red-highlighted C) instead of inheriting
# the bold style of column A from the row above.
$ws.Range("B2:C2").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Remove any formatting/content accidentally carried into column A of the
# newly inserted row - the target has no A13 cell at all.
$ws.Cells.Item(13, 1).Clear()

# Populate the newly inserted row with the "Docentes responsáveis:" data
# (previously duplicated/misplaced at row 18 under "Método:").
$docentesTxt = "8767640 - Eduardo Ferro dos Santos"
$ws.Cells.Item(13, 2).Value = $docentesTxt
$ws.Cells.Item(13, 3).Value = $docentesTxt

# 3) Update B/C text on the shifted rows (14-22); the A-column labels are
#    unchanged by the shift, only the B/C contents differ from before.
$programaResumidoTxt = "Introdução aos softwares de desenho técnico. Primitivas gráficas e operações de manipulação. Desenhos 2D. Desenhos 3D. Importação e exportação de dados. Aplicações especiais. Montagem de conjuntos. Visualização de movimentos e interferências."
$ws.Cells.Item(14, 2).Value = $programaResumidoTxt
$ws.Cells.Item(14, 3).Value = $programaResumidoTxt

$programaTxt = "Introdução ao hardware e periféricos gráficos; Estrutura genérica de softwares de desenho técnico; Primitivas gráficas e operações de manipulação de elementos gráficos; Execução de desenhos com softwares de desenho técnico: Desenhos 2D, Desenhos 3D, Vistas a partir de modelos 3D, Importação e exportação de dados, Uso de Bibliotecas gráficas, Esboços e modelamento paramétrico, Aplicações especiais, Montagem de conjuntos, Visualização de movimentos e interferências."
$ws.Cells.Item(16, 2).Value = $programaTxt
$ws.Cells.Item(16, 3).Value = $programaTxt

$metodoTxt = "Aulas expositivas e práticas."
$ws.Cells.Item(19, 2).Value = $metodoTxt
$ws.Cells.Item(19, 3).Value = $metodoTxt

$criterioTxt = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"
$ws.Cells.Item(20, 2).Value = $criterioTxt
$ws.Cells.Item(20, 3).Value = $criterioTxt

$normaTxt = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
$ws.Cells.Item(21, 2).Value = $normaTxt
$ws.Cells.Item(21, 3).Value = $normaTxt

$bibliografiaTxt = "Bibliografia principal:- Manual do software a ser adotado. Tutorias do software a ser adotado.Bibliografia complementar:- Normas Brasileiras Aplicadas ao Desenho Técnico.- RIBEIRO, A. C.; PERES, M. P.; IZIDORO, N. Curso de Desenho Técnico e AutoCAD, Pearson, 2013.- SILVA, A., RIBEIRO, C.T., DIAS, J. e SOUSA, L. Desenho Técnico Moderno, LTC, 2006.- LEAKE, J. e BORGERSON, J. Manual de Desenho Técnico para Engenharia - Desenho, Modelagem e Visualização. LTC, 2010."
$ws.Cells.Item(22, 2).Value = $bibliografiaTxt
$ws.Cells.Item(22, 3).Value = $bibliografiaTxt
